$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "32.00", "1.00") are stored as exact text rather than being
# reinterpreted as numbers. Style is reset back to Normal afterwards so
# the cell styling matches the original (no explicit style index).
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range('D2').Value = '61.392.29'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '2.884.18'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '567.12'
$ws.Range('E5').Value = '  -4.28%  '
$ws.Range('D6').Value = '143.74'
$ws.Range('E6').Value = '  -2.44%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '2.883.39'
$ws.Range('D10').Value = '6.94'
$ws.Range('E10').Value = '  -4.82%  '
$ws.Range('E11').Value = '  -1.38%  '
$ws.Range('E12').Value = '  -1.83%  '
$ws.Range('E13').Value = '  -0.50%  '
$ws.Range('D14').Value = '32.00'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '3.354.38'
$ws.Range('E16').Value = '  -2.16%  '
$ws.Range('D17').Value = '61.381.04'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = '2.907.44'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('E20').Value = '  -1.77%  '
$ws.Range('D21').Value = '13.06'
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').Value = '0.656'
$ws.Range('E22').Value = '  -1.16%  '
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('D24').Value = '79.32'
$ws.Range('E24').Value = '  -1.84%  '
$ws.Range('D25').Value = '11.90'
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').Value = '10.02'
$ws.Range('E26').Value = '  -9.68%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('E28').Value = '  -5.38%  '
$ws.Range('E29').Value = '  +4.44%  '
$ws.Range('E31').Value = '  -3.84%  '
$ws.Range('E32').Value = '  -7.46%  '
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('D34').Value = '0.107'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('D35').Value = '25.41'
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('D36').Value = '0.959'
$ws.Range('E36').Value = '  -3.20%  '
$ws.Range('D37').Value = '5.38'
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('D38').Value = '48.87'
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('E39').Value = '  -7.41%  '
$ws.Range('E40').Value = '  -4.20%  '
$ws.Range('D41').Value = '8.21'
$ws.Range('E41').Value = '  -2.85%  '
$ws.Range('E42').Value = '  -2.74%  '
$ws.Range('D43').Value = '38.97'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('D45').Value = '2.709.15'
$ws.Range('D46').Value = '133.10'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('D47').Value = '0.0336'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').Value = '342.59'
$ws.Range('E48').Value = '  -4.67%  '
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').Value = '21.57'
$ws.Range('E51').Value = '  -4.40%  '

$colD.Style = "Normal"
